$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-44: update Price (D) and Volume(1h) (E) values
# D column is forced to Text format first so numeric-looking strings
# (e.g. "1.000", "0.9998") are preserved exactly as text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.099.62"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.53"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.06"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6846"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3010"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07429"
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.04"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07644"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.46"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.048"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6814"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.36"
$ws.Range("E15").Value = "  -6.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.167"
$ws.Range("E16").Value = "  -6.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.095.19"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008130"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.080.06"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.74"
$ws.Range("E20").Value = "  -5.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.377"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.97"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1449"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.744"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.515"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.273"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.129"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05224"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7578"
$ws.Range("E34").Value = "  -4.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.844"
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.132"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.301.26"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01833"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.722"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9312"
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.839"
$ws.Range("E42").Value = "  -4.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.73"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.01%  "

# Rows 45-46: RocketPoolETH and BabyDogeCoin swap list positions, with updated price/volume
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000130"
$ws.Range("E45").Value = "  +5.13%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.985.15"
$ws.Range("E46").Value = "  -0.29%  "

# Rows 47-51: update Price (D) and Volume(1h) (E) values
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5197"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.72"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.488"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.766"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05944"
